$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.947.41'
$ws.Range('E2').Value = '  -1.18%  '

$ws.Range('D3').Value = '3.413.93'
$ws.Range('E3').Value = '  -1.50%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '404.54'
$ws.Range('E5').Value = '  -0.91%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.18'
$ws.Range('E6').Value = '  +1.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.590'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.688'
$ws.Range('E9').Value = '  -1.54%  '

$ws.Range('E10').Value = '  -3.80%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.80'
$ws.Range('E11').Value = '  -3.47%  '

$ws.Range('E12').Value = '  -1.21%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.41'
$ws.Range('E13').Value = '  -4.53%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.77'
$ws.Range('E14').Value = '  -1.76%  '

$ws.Range('D15').Value = '3.538.55'
$ws.Range('E15').Value = '  +2.13%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '11.68'
$ws.Range('E16').Value = '  +7.92%  '

$ws.Range('D17').Value = '61.979.42'
$ws.Range('E17').Value = '  -1.25%  '

$ws.Range('E18').Value = '  -3.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000141'
$ws.Range('E19').Value = '  -0.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.16'
$ws.Range('E20').Value = '  -5.40%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '83.56'
$ws.Range('E21').Value = '  +0.80%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '311.88'
$ws.Range('E22').Value = '  -0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.74'
$ws.Range('E23').Value = '  -3.63%  '

$ws.Range('E24').Value = '  -1.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.76'
$ws.Range('E25').Value = '  +8.95%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.58'
$ws.Range('E26').Value = '  -2.97%  '

$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.88'
$ws.Range('E27').Value = '  +3.73%  '

$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.04'
$ws.Range('E28').Value = '  -2.30%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.74'
$ws.Range('E29').Value = '  +5.22%  '

$ws.Range('E30').Value = '  -2.58%  '

$ws.Range('E31').Value = '  -2.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '42.47'
$ws.Range('E32').Value = '  -3.73%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('E34').Value = '  -4.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0481'
$ws.Range('E35').Value = '  -2.83%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.48'
$ws.Range('E36').Value = '  -2.18%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.329'
$ws.Range('E38').Value = '  +13.72%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.37'
$ws.Range('E39').Value = '  -6.34%  '

$ws.Range('E40').Value = '  -3.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.41'
$ws.Range('E41').Value = '  +2.06%  '

$ws.Range('E42').Value = '  -1.35%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.96'
$ws.Range('E43').Value = '  -1.14%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.96'
$ws.Range('E44').Value = '  -0.36%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.55'
$ws.Range('E45').Value = '  -6.24%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.22'
$ws.Range('E46').Value = '  -1.06%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.22'
$ws.Range('E47').Value = '  -4.01%  '

$ws.Range('D48').Value = '2.106.13'
$ws.Range('E48').Value = '  -3.52%  '

$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.82'
$ws.Range('E49').Value = '  +25.25%  '

$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.35'
$ws.Range('E50').Value = '  -3.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.94'
$ws.Range('E51').Value = '  +3.06%  '
